$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed crypto price / 1h-volume snapshot figures.
# All values in columns D and E are stored as plain text in the source sheet
# (mixed "."-grouped prices, and percentages padded with spaces), so cells whose
# new content looks like a plain decimal number (e.g. "200.35") must be forced to
# stay text - otherwise Excel auto-converts them to a numeric value on assignment.

$ws.Range("D2").Value = "76.984.65"
$ws.Range("E2").Value = "  +0.88%  "
$ws.Range("D3").Value = "2.965.17"
$ws.Range("E3").Value = "  +3.25%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "200.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "596.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.552"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.197"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.96%  "
$ws.Range("D10").Value = "2.963.54"
$ws.Range("E10").Value = "  +3.24%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.443"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +13.03%  "
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "3.510.21"
$ws.Range("E13").Value = "  +3.37%  "
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "76.702.92"
$ws.Range("E15").Value = "  +0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.16%  "
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").Value = "2.955.92"
$ws.Range("E18").Value = "  +3.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.48"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.61%  "
$ws.Range("E20").Value = "  -3.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "373.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.08%  "
$ws.Range("E22").Value = "  +4.87%  "
$ws.Range("E23").Value = "  -3.64%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "72.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.85%  "
$ws.Range("D25").Value = "3.116.28"
$ws.Range("E25").Value = "  +2.88%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.27"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.72"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("E29").Value = "  +3.37%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.09%  "
$ws.Range("E32").Value = "  -1.50%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "499.35"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.94%  "
$ws.Range("E34").Value = "  +1.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.114"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +24.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.22"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.61%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.398"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "20.19"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  +1.37%  "
$ws.Range("E41").Value = "  -6.24%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "180.50"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.27%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "4.94"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.31%  "
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("E47").Value = "  -2.90%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.591"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.40%  "
$ws.Range("E49").Value = "  +4.25%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.31"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "22.67"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.06%  "
